$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (and update the title text used elsewhere)
$ws.Name = "Through 2022-09-08"

# Update the label for the September row
$ws.Range("A10").Value = "September (through 09-08)"

# Update September row values (row 10)
$ws.Range("B10").Value = 8
$ws.Range("D10").Value = 19
$ws.Range("E10").Value = 11
$ws.Range("F10").Value = 20
$ws.Range("G10").Value = 27
$ws.Range("H10").Value = 32
$ws.Range("I10").Value = 38

# Update Total row values (row 11)
$ws.Range("B11").Value = 202
$ws.Range("D11").Value = 570
$ws.Range("E11").Value = 501
$ws.Range("F11").Value = 369
$ws.Range("G11").Value = 811
$ws.Range("H11").Value = 1102
$ws.Range("I11").Value = 1175
